$d = $word.ActiveDocument

# Common paragraph-property blocks reused below.
$pPr0 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr>'
$pPr7 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/></w:pPr>'
$pPr7Color = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:color w:val="00B050"/></w:rPr></w:pPr>'

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
$pkgClose = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1) "Become familiar with Sifteo API (all)" -> split run, wrap "Sifteo" in
#    proofErr spell-check markers.
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$inner4 = $pPr7Color +
    '<w:r w:rsidRPr="005D2912"><w:rPr><w:color w:val="00B050"/></w:rPr><w:t xml:space="preserve">Become familiar with </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>Sifteo</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t xml:space="preserve"> API (all)</w:t></w:r>'
$p4.Range.InsertXML($pkgOpen + $inner4 + $pkgClose)

# ---------------------------------------------------------------------------
# 2) "Conclude Kivy investigation (Kurtis/Alex)" -> split run, wrap "Kivy".
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$inner5 = $pPr7Color +
    '<w:r w:rsidRPr="005D2912"><w:rPr><w:color w:val="00B050"/></w:rPr><w:t xml:space="preserve">Conclude </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>Kivy</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t xml:space="preserve"> investigation (Kurtis/Alex)</w:t></w:r>'
$p5.Range.InsertXML($pkgOpen + $inner5 + $pkgClose)

# ---------------------------------------------------------------------------
# 3) Move the _GoBack bookmark from the "Implement cube actions in UI"
#    paragraph to the end of "Action items for this week", and retitle the
#    action-item list (each item's text shifts down to make room for a new
#    "Implement Cube.FillScreen()" item at the end, and item 1 becomes a
#    brand-new "Refine cube click-and-drag" item).
# ---------------------------------------------------------------------------
$p11 = $d.Paragraphs(11)
$inner11 = $pPr0 + '<w:r><w:t>Action items for this week</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$p11.Range.InsertXML($pkgOpen + $inner11 + $pkgClose)

$p12 = $d.Paragraphs(12)
$inner12 = $pPr7 + '<w:r><w:t>Refine cube click-and-drag</w:t></w:r>'
$p12.Range.InsertXML($pkgOpen + $inner12 + $pkgClose)

$p13 = $d.Paragraphs(13)
$inner13 = $pPr7 + '<w:r><w:t>Implement cube actions in UI</w:t></w:r>'
$p13.Range.InsertXML($pkgOpen + $inner13 + $pkgClose)

$p14 = $d.Paragraphs(14)
$inner14 = $pPr7 +
    '<w:r><w:t xml:space="preserve">Implement </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Cube.Paint</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>()</w:t></w:r>'
$p14.Range.InsertXML($pkgOpen + $inner14 + $pkgClose)

$p15 = $d.Paragraphs(15)
$inner15 = $pPr7 +
    '<w:r><w:t xml:space="preserve">Implement </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>BaseApp.Setup</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>()</w:t></w:r>'
$p15.Range.InsertXML($pkgOpen + $inner15 + $pkgClose)

$p16 = $d.Paragraphs(16)
$inner16 = $pPr7 +
    '<w:r><w:t xml:space="preserve">Implement </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>BaseApp.Tick</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>()</w:t></w:r>'
$p16.Range.InsertXML($pkgOpen + $inner16 + $pkgClose)

$p17 = $d.Paragraphs(17)
$inner17 = $pPr7 +
    '<w:r><w:t xml:space="preserve">Implement </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Cube.FillRect</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>()</w:t></w:r>'
$p17.Range.InsertXML($pkgOpen + $inner17 + $pkgClose)

# Insert a brand-new paragraph right after (old) item 17, inheriting its
# list formatting, to host the new "Implement Cube.FillScreen()" line.
$p17.Range.InsertParagraphAfter()

$p18 = $d.Paragraphs(18)
$inner18 = $pPr7 +
    '<w:r><w:t xml:space="preserve">Implement </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Cube.FillScreen</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>()</w:t></w:r>'
$p18.Range.InsertXML($pkgOpen + $inner18 + $pkgClose)

Write-Host "Done"
